$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "<and is tential hundrate attention, wheres that the advanced in betten bearshian, wheres that the advanced in betth advanced.>"
$ws.Range("C2").Value = 59.88277483587631

$ws.Range("B3").Value = "<three months mothers who his mother to his mother to his mothers who his mother to his mother work strier to his mother to his mother to his mother work.>"
$ws.Range("C3").Value = 60.2910097338486

$ws.Range("B4").Value = "<in nenep thousand perion the nenep thousand rickery and perion opwards.>"
$ws.Range("C4").Value = 61.53344785272301

$ws.Range("B5").Value = "<and waited in the track in the track and waited in the track in the track in the corner.>"
$ws.Range("C5").Value = 60.38171760953804

$ws.Range("B6").Value = "<and to go into the contry on the business of the contry on the business of the contry on the business of the contry one.>"
$ws.Range("C6").Value = 60.45639681782397

$ws.Range("B7").Value = "<he deber to get curtain to get aver to get aver to get aver to get curtain to gerther rods for and ronted aborant.>"
$ws.Range("C7").Value = 59.44948501942512

$ws.Range("B8").Value = "<wat the held as the opport the held at the held a hellas the held ar.>"
$ws.Range("C8").Value = 58.90765918201383

$ws.Range("C9").Value = 62.08122809870977

$ws.Range("B10").Value = "<in the marching the conviction asses he evidence of march the conviction as to the conviction.>"
$ws.Range("C10").Value = 60.21943983566331

$ws.Range("B11").Value = "<to the bitter distruction of all and every distruction of all and every distruction of all and every distruction>"
$ws.Range("C11").Value = 63.10557367411857

$ws.Range("B12").Value = "<but they never perpoplished they neves and acomplished they nethems and acomplished the man at on eless of manate the man at on easfort of mannes.>"
$ws.Range("C12").Value = 60.88751564398321

$ws.Range("B13").Value = "<in the case all licked in flicted murders of unishment in flicted murders olly,>"
$ws.Range("C13").Value = 62.01012305716

$ws.Range("B14").Value = "<three, condemnement to death and three, condemnement to death and three, condemnement to death an executing executing executing executing.>"
$ws.Range("C14").Value = 60.98431102175819

$ws.Range("B15").Value = "<a great states thing dust of a changemage oswald states thing duse him dustates,>"
$ws.Range("C15").Value = 58.8825103220702

$ws.Range("B16").Value = "<he was stear people good works the atter p.m. acted the low could who could who chrominan all good works.>"
$ws.Range("C16").Value = 63.75228872147132

$ws.Range("B17").Value = "<main doorsisty four,>"
$ws.Range("C17").Value = 61.62584077095502

$ws.Range("B18").Value = "<ows starle mother explose starle mother explose starley>"
$ws.Range("C18").Value = 62.39030428536886

$ws.Range("B19").Value = "<he was could hamped to a part in a day went roker fren a day went roker fred a hamped to a part in a day was coant roker freen at roker fridend.>"
$ws.Range("C19").Value = 64.48089963139564

$ws.Range("B20").Value = "<and interess should be made for the atterate of the should be made for the atterate a vatery,>"
$ws.Range("C20").Value = 65.99867469216844

$ws.Range("B21").Value = "<some of these founds of the wroims of the wroims of the wroims of the wrooms of the walls the walls the walls the walls.>"
$ws.Range("C21").Value = 61.79385642216577

$ws.Range("B22").Value = "<oswald neing in inchincies five for a years all was five for ander and for years,>"
$ws.Range("C22").Value = 63.78536900073647

$ws.Range("B23").Value = "<the service proform and vansuations proform and vans proformations.>"
$ws.Range("C23").Value = 66.07132560035792

$ws.Range("B24").Value = "<however, a man tlange with plain lange with plain lange lange lange with plain lange lain lange down in tlain.>"
$ws.Range("C24").Value = 59.96065374062437

$ws.Range("B25").Value = "<the game however, walking however, as met by her was near was met by her was near was met by her was near was near was near was near walond.>"
$ws.Range("C25").Value = 61.37294568965083

$ws.Range("B26").Value = "<when he had the words react he had the words react he had the words react smistured smistured smistewaled#>"
$ws.Range("C26").Value = 61.49304647260347

$ws.Range("B27").Value = "<as well three, dated mated mated mated mated mated mated mated mated mated mated mate>"
$ws.Range("C27").Value = 57.65014137385923

$ws.Range("B28").Value = "<five six points points points points points.>"
$ws.Range("C28").Value = 60.19432755083982

$ws.Range("B29").Value = "<oswald#s rivolver,>"
$ws.Range("C29").Value = 59.59015377982527

$ws.Range("B30").Value = "<lasisor loorned thim against is renst sis relations warned.>"
$ws.Range("C30").Value = 65.73334965659353

$ws.Range("B31").Value = "<this committee and versially upon the conmittee, and wares, and worse and upon the conmittee, and warest commutioned strongly upon the condit>"
$ws.Range("C31").Value = 64.11428951912866

$ws.Range("B32").Value = "<in regarding the commission regarding the and motorcade body by the commission regarding the incording>"
$ws.Range("C32").Value = 59.41950983700274

$ws.Range("B33").Value = "<two has been the new orleans member in the new orleans member should car.>"
$ws.Range("C33").Value = 60.9158653852691

$ws.Range("B34").Value = "<was then the president#s car at the speated who dinct with regate who dinct with regar at the presidents car.>"
$ws.Range("C34").Value = 63.86599029640286

$ws.Range("B35").Value = "<his appearable for for for ferman sented a prole for for ferman sented a prole for for for for ferman sentement in new marinst from in new marice prole for the brip.>"
$ws.Range("C35").Value = 59.12062429245024

$ws.Range("B36").Value = "<advoisive that oswald was on no is the bure of that oswald#s.>"
$ws.Range("C36").Value = 59.46977328978135

$ws.Range("B37").Value = "<in the second serit to trand, quote, fandual paper to a real paper to a real paper to a real paper to a real paper to arrit to trand, quote.>"
$ws.Range("C37").Value = 64.65723548154703

$ws.Range("B38").Value = "<calf the suppliet thirty polight half of the forgoing ritty polight polight polight polight polight polings.>"
$ws.Range("C38").Value = 63.08988765490851

$ws.Range("B39").Value = "<to have though have though have though have though have though have though have though have though have a tenses.>"
$ws.Range("C39").Value = 61.03478338354731

$ws.Range("B40").Value = "<and told him a similar similar story told him a story told himse.>"
$ws.Range("C40").Value = 62.65136469847736

$ws.Range("B41").Value = "<which her end quote, which him and ven ven, end quote.>"
$ws.Range("C41").Value = 65.79040011990128

$ws.Range("B42").Value = "<they talked a days were hump days were hump days were hump days were hump days were hum.>"
$ws.Range("C42").Value = 59.09579034365078
